$wb = $excel.ActiveWorkbook

# --- Sheet ALC: scheduled market-price refresh ---
$ws = $wb.Sheets("ALC")
$ws.Range("H80").Value = 4259.7
$ws.Range("I80").Value = 4116.5
$ws.Range("J80").Value = 4474.5
$ws.Range("K80").Value = 12349.5
$ws.Range("L80").Value = 13423.5
$ws.Range("M80").Value = -11351.5
$ws.Range("N80").Value = -15419.5
$ws.Range("H83").Value = 4259.7
$ws.Range("I83").Value = 4116.5
$ws.Range("J83").Value = 4474.5
$ws.Range("K83").Value = 37048.5
$ws.Range("L83").Value = 40270.5
$ws.Range("M83").Value = -32056.5
$ws.Range("N83").Value = -50254.5
$ws.Range("H86").Value = 6571.2856
$ws.Range("I86").Value = 5799.8
$ws.Range("J86").Value = 8500
$ws.Range("K86").Value = 5799.8
$ws.Range("L86").Value = 8500
$ws.Range("M86").Value = -4676.8
$ws.Range("N86").Value = -10746
$ws.Range("H89").Value = 6571.2856
$ws.Range("I89").Value = 5799.8
$ws.Range("J89").Value = 8500
$ws.Range("K89").Value = 28999
$ws.Range("L89").Value = 42500
$ws.Range("M89").Value = -23383
$ws.Range("N89").Value = -53732
$ws.Range("H111").Value = 1700.5
$ws.Range("I111").Value = 774.7778
$ws.Range("K111").Value = 2324.3334
$ws.Range("M111").Value = 742.6666
$ws.Range("H116").Value = 9900
$ws.Range("I116").Value = 9900
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 9900
$ws.Range("L116").Value = 0
$ws.Range("N116").Value = -6458
$ws.Range("H138").Value = 4050.0557
$ws.Range("J138").Value = 5534.3335
$ws.Range("L138").Value = 16603.0005
$ws.Range("N138").Value = -26883.0005
$ws.Range("H141").Value = 533.3333
$ws.Range("I141").Value = 600
$ws.Range("J141").Value = 400
$ws.Range("K141").Value = 1800
$ws.Range("L141").Value = 1200
$ws.Range("M141").Value = 3380
$ws.Range("N141").Value = -11560

# --- Sheet ARM: scheduled market-price refresh ---
$ws = $wb.Sheets("ARM")
$ws.Range("H61").Value = 1249.5
$ws.Range("I61").Value = 1288.3334
$ws.Range("K61").Value = 1288.3334
$ws.Range("M61").Value = -1076.3334
$ws.Range("H88").Value = 3001
$ws.Range("I88").Value = 1403
$ws.Range("J88").Value = 3800
$ws.Range("K88").Value = 1403
$ws.Range("L88").Value = 3800
$ws.Range("M88").Value = -997
$ws.Range("N88").Value = -4612
$ws.Range("H91").Value = 3001
$ws.Range("I91").Value = 1403
$ws.Range("J91").Value = 3800
$ws.Range("K91").Value = 1403
$ws.Range("L91").Value = 3800
$ws.Range("M91").Value = 1
$ws.Range("N91").Value = -6608
$ws.Range("H102").Value = 4087.5715
$ws.Range("I102").Value = 3935.5
$ws.Range("K102").Value = 3935.5
$ws.Range("M102").Value = -2313.5
$ws.Range("H136").Value = 1249.5
$ws.Range("I136").Value = 1288.3334
$ws.Range("K136").Value = 3865.0002
$ws.Range("M136").Value = -1315.0002

# --- Sheet BSM: scheduled market-price refresh ---
$ws = $wb.Sheets("BSM")
$ws.Range("H86").Value = 5798.1665
$ws.Range("I86").Value = 2563
$ws.Range("J86").Value = 9033.333000000001
$ws.Range("K86").Value = 2563
$ws.Range("L86").Value = 9033.333000000001
$ws.Range("M86").Value = -1440
$ws.Range("N86").Value = -11279.333
$ws.Range("H89").Value = 5798.1665
$ws.Range("I89").Value = 2563
$ws.Range("J89").Value = 9033.333000000001
$ws.Range("K89").Value = 12815
$ws.Range("L89").Value = 45166.665
$ws.Range("M89").Value = -7199
$ws.Range("N89").Value = -56398.665
$ws.Range("H105").Value = 4999.6665
$ws.Range("I105").Value = 4999.6665
$ws.Range("K105").Value = 4999.6665
$ws.Range("M105").Value = -3252.6665
$ws.Range("H134").Value = 3858.4
$ws.Range("I134").Value = 3835.5
$ws.Range("K134").Value = 11506.5
$ws.Range("M134").Value = -8971.5

# --- Sheet CRP: scheduled market-price refresh ---
$ws = $wb.Sheets("CRP")
$ws.Range("H31").Value = 2189.7
$ws.Range("J31").Value = 2319.8
$ws.Range("L31").Value = 2319.8
$ws.Range("N31").Value = -2909.8
$ws.Range("H34").Value = 2189.7
$ws.Range("J34").Value = 2319.8
$ws.Range("L34").Value = 2319.8
$ws.Range("N34").Value = -2723.8
$ws.Range("H134").Value = 6547
$ws.Range("I134").Value = 7156.1665
$ws.Range("K134").Value = 21468.4995
$ws.Range("M134").Value = -18933.4995

# --- Sheet CUL: scheduled market-price refresh ---
$ws = $wb.Sheets("CUL")
$ws.Range("H4").Value = 3809.375
$ws.Range("I4").Value = 4343.5713
$ws.Range("K4").Value = 13030.7139
$ws.Range("M4").Value = -12918.7139

# --- Sheet GSM: scheduled market-price refresh ---
$ws = $wb.Sheets("GSM")
$ws.Range("H80").Value = 7499.8335
$ws.Range("J80").Value = 9749.75
$ws.Range("L80").Value = 9749.75
$ws.Range("N80").Value = -11745.75
$ws.Range("H83").Value = 7499.8335
$ws.Range("J83").Value = 9749.75
$ws.Range("L83").Value = 48748.75
$ws.Range("N83").Value = -58732.75
$ws.Range("H136").Value = 41404.332
$ws.Range("J136").Value = 41404.332
$ws.Range("L136").Value = 124212.996
$ws.Range("N136").Value = -129312.996

# --- Sheet LTW: scheduled market-price refresh ---
$ws = $wb.Sheets("LTW")
$ws.Range("H46").Value = 2493.25
$ws.Range("J46").Value = 5083
$ws.Range("L46").Value = 5083
$ws.Range("N46").Value = -5459
$ws.Range("H82").Value = 2688.125
$ws.Range("I82").Value = 2853.3333
$ws.Range("K82").Value = 2853.3333
$ws.Range("M82").Value = -2492.3333
$ws.Range("H85").Value = 2688.125
$ws.Range("I85").Value = 2853.3333
$ws.Range("K85").Value = 2853.3333
$ws.Range("M85").Value = -1605.3333

# --- Sheet WVR: scheduled market-price refresh ---
$ws = $wb.Sheets("WVR")
$ws.Range("H62").Value = 6235
$ws.Range("I62").Value = 2625
$ws.Range("K62").Value = 2625
$ws.Range("M62").Value = -2001
$ws.Range("H65").Value = 6235
$ws.Range("I65").Value = 2625
$ws.Range("K65").Value = 13125
$ws.Range("M65").Value = -10005
$ws.Range("H81").Value = 4167.6
$ws.Range("I81").Value = 1852.8889
$ws.Range("J81").Value = 25000
$ws.Range("K81").Value = 3705.7778
$ws.Range("L81").Value = 50000
$ws.Range("M81").Value = -2644.7778
$ws.Range("N81").Value = -52122
$ws.Range("H84").Value = 4167.6
$ws.Range("I84").Value = 1852.8889
$ws.Range("J84").Value = 25000
$ws.Range("K84").Value = 18528.889
$ws.Range("L84").Value = 250000
$ws.Range("M84").Value = -13224.889
$ws.Range("N84").Value = -260608
